$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = "BR"
$ws.Range("C32").Value = 200
$ws.Range("D32").Value = "Fecha final debe ser mayor a la fecha inicial"

$ws.Rows(4).Select() | Out-Null
